$d = $word.ActiveDocument

# --- 1. Turn the "[24/8/20]" log-entry paragraph into a Heading3 ---
$datePara = $d.Paragraphs.Item(3)
$datePara.Style = "Heading3"

# Bookmark the date text itself (exclude the trailing paragraph mark)
$dateRange = $d.Range($datePara.Range.Start, $datePara.Range.End - 1)
$d.Bookmarks.Add("section", $dateRange)

# --- 2. The paragraph right after it switches from BodyText to FirstParagraph ---
$firstEntryPara = $d.Paragraphs.Item(4)
$firstEntryPara.Style = "FirstParagraph"

# --- 3. Append a new log entry ("[24/8/27]") at the end of the document ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newHeading = $d.Paragraphs.Last
$newHeading.Style = "Heading3"
$newHeading.Range.Text = "[24/8/27]"

# Bookmark the new heading's date text (exclude the trailing paragraph mark)
$newHeadingRange = $d.Range($newHeading.Range.Start, $newHeading.Range.End - 1)
$d.Bookmarks.Add("section-1", $newHeadingRange)

$newHeading.Range.InsertParagraphAfter()

$newBody = $d.Paragraphs.Last
$newBody.Style = "FirstParagraph"
$newBody.Range.Text = "Idag har vi gjort klart tornet. Vi tacklade också problemet av stabilitet genom att sätta på en bas på ostabila sidan av tornet. Nu verkar vårt torn kunna hålla väldigt mycket vikt och är ändå rätt högt."
